$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-12-17 Tuesday" "2024-12-18 Wednesday"

Replace-Text "416×7=" "867×9="
Replace-Text "937×9=" "268×7="
Replace-Text "835×4=" "428×2="
Replace-Text "450×5=" "606×2="
Replace-Text "245×2=" "760×2="

Replace-Text "178×3=" "671×6="
Replace-Text "443×7=" "361×5="
Replace-Text "234×4=" "557×2="
Replace-Text "894×3=" "839×4="
Replace-Text "602×5=" "857×6="

Replace-Text "807×8=" "953×5="
Replace-Text "636×4=" "693×6="
Replace-Text "606×9=" "528×4="
Replace-Text "529×4=" "644×4="
Replace-Text "705×2=" "672×7="

Replace-Text "426×5=" "266×4="
Replace-Text "903×5=" "743×8="
Replace-Text "311×9=" "796×9="
Replace-Text "200×4=" "338×9="
Replace-Text "828×2=" "333×5="

Replace-Text "172×9=" "566×4="
Replace-Text "122×8=" "549×4="
Replace-Text "846×8=" "602×2="
Replace-Text "560×3=" "472×7="
Replace-Text "703×8=" "497×8="
